# C5-PowerPoint.pptx edit
# 1) Re-style the "Sources of finance" table on slide 6 with the new table style GUID.
# 2) Swap the deck's applied colour theme from "Integral" to the classic "Office" palette
#    (font/effect schemes are identical between the two themes - only the 12 theme colours
#    differ), mirroring the theme1.xml <-> theme2.xml colour swap from the source commit.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 ------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{6AD99924-EFC2-4E30-89E4-D1199579B2B2}")

# --- 2. Theme colours ----------------------------------------------------------
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

$colors.Colors(1).RGB  = 0         # dk1      000000
$colors.Colors(2).RGB  = 16777215  # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388   # dk2      44546A
$colors.Colors(4).RGB  = 15132391  # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939  # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501   # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845  # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407     # accent4  FFC000
$colors.Colors(9).RGB  = 12874308  # accent5  4472C4
$colors.Colors(10).RGB = 4697456   # accent6  70AD47
$colors.Colors(11).RGB = 12673797  # hlink    0563C1
$colors.Colors(12).RGB = 7491477   # folHlink 954F72
